# Add labels for newly added roles
#
# The "[Roles]" table on the active sheet has columns:
#   A = role code, B = RoleName, C = label
# A handful of rows were missing their label value in column C. Fill them
# in with the same text already used for that role in columns A/B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C28").Value = $ws.Range("A28").Value2
$ws.Range("C30").Value = $ws.Range("A30").Value2
$ws.Range("C31").Value = $ws.Range("A31").Value2
$ws.Range("C32").Value = $ws.Range("A32").Value2

# Row 33 never had any value at all in column C, so the new cell would
# otherwise inherit the plain column default formatting instead of
# matching the rest of the row. Copy the formatting from a neighbouring
# cell in the same row first, then fill in the label.
[void]$ws.Range("D33").Copy()
[void]$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("C33").Value = $ws.Range("A33").Value2

# Leave the selection where the editor ended up after making the change.
[void]$ws.Range("C35").Select()
